$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: strike-through the "Дочитать до конца GoF book" paragraph
#           (paragraph mark + every run gets <w:strike/>)
# ---------------------------------------------------------------------------
$strikeRange = $d.Content
$strikeRange.Find.Execute("Дочитать до конца GoF book") | Out-Null
$strikePara = $strikeRange.Paragraphs(1)
$strikePara.Range.Font.StrikeThrough = $true

# ---------------------------------------------------------------------------
# Change 2: insert a new yellow-highlighted bullet just before the
#           "Концепция нескольких мастер-таблиц..." paragraph
# ---------------------------------------------------------------------------
$targetRange = $d.Content
$targetRange.Find.Execute("Концепция нескольких") | Out-Null
$targetPara = $targetRange.Paragraphs(1)
$newIndex = $targetPara.Index
$targetPara.Range.InsertParagraphBefore() | Out-Null

# Re-locate the freshly created (still empty) paragraph by its numeric index
# and fill it in with the exact WordprocessingML the commit introduces
# (including the spell-check markers around "MasterTable").
$newPara = $d.Paragraphs($newIndex)

$xml = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="a3"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr><w:rPr><w:highlight w:val="yellow"/></w:rPr></w:pPr><w:r><w:rPr><w:highlight w:val="yellow"/></w:rPr><w:t xml:space="preserve">Объединять ячейки таблицы с данными, относящиеся к одной записи </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:highlight w:val="yellow"/></w:rPr><w:t>MasterTable</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:highlight w:val="yellow"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$newPara.Range.InsertXML($xml) | Out-Null
